# Update the dSF column (F) values for several rows in Sheet1.
# This reflects a "repull data, push all data, mean calculation" update
# where the final dS (dSF) values were recalculated and now differ
# from the initial dS0 values for a subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    3  = 8
    4  = 7
    5  = -1
    13 = -2
    16 = -1
    17 = -1
    18 = 2
    21 = 0
    24 = 5
    27 = -6
    30 = -3
    32 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
